$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new "Period" column before the existing "Value" column (B),
# shifting Value/As Of/Notes one column to the right.
$ws.Columns("B").Insert()

# Header for the new column
$ws.Range("B1").Value = "Period"

# Populate the Period values: quarterly KPIs except "Num of Stores", which is Annual
$ws.Range("B2").Value = "Quarter"
$ws.Range("B3").Value = "Quarter"
$ws.Range("B4").Value = "Quarter"
$ws.Range("B5").Value = "Annual"

# Match the new column's width to column A
$ws.Columns("B").ColumnWidth = $ws.Columns("A").ColumnWidth

# Leave the selection where the user would naturally end up after typing the last entry
$ws.Range("B6").Select() | Out-Null
